$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Estado C1" (column M) values for each test case row from
# "PENDIENTE" to the actual outcome of the test run.
$ws.Range("M3").Value = "Exitoso"
$ws.Range("M4").Value = "Exitoso"
$ws.Range("M5").Value = "Exitoso"
$ws.Range("M6").Value = "No Exitoso"
$ws.Range("M7").Value = "Exitoso"
$ws.Range("M8").Value = "No Exitoso"
$ws.Range("M9").Value = "Exitoso"

# Leave the merged "EVIDENCIA" header cell selected, matching where the
# user's cursor ended up after updating the states.
$ws.Range("M1:P1").Select()
